$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 175
$ws.Range("J2").Value = 694
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 184
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 122
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 15
$ws.Range("S2").Value = 88
$ws.Range("T2").Value = 118
$ws.Range("U2").Value = 6
$ws.Range("V2").Value = 1025
$ws.Range("X2").Value = 999
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 7
